$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Checklist" to "Session"
$ws.Name = "Session"

# Update column E ("Type") values from "Selection" to "Scan" for data rows 2-52
for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    if ($cell.Value() -eq "Selection") {
        $cell.Value = "Scan"
    }
}
